$d = $word.ActiveDocument

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $ns + '><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Set-ParagraphXml($label, $paragraphRange, $bodyInnerXml) {
    if ($null -eq $paragraphRange) {
        Write-Host "ERROR: target paragraph range not found for $label"
        return
    }
    $xml = $pkgHeader + $bodyInnerXml + $pkgFooter
    $paragraphRange.InsertXML($xml) | Out-Null
    Write-Host "OK: $label"
}

function Find-ParagraphContaining($needle) {
    $paras = $d.Paragraphs
    $count = $paras.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        if ($p.Range.Text -like "*$needle*") {
            return $p.Range
        }
    }
    return $null
}

# Some phrases (e.g. the ACM Code of Ethics heading) occur more than once in
# the document (once in the intro table of contents, once as the big-font
# section heading, and possibly once more in body text). Disambiguate by
# also requiring a specific font size on the matching paragraph.
function Find-ParagraphContainingWithFontSize($needle, $fontSize) {
    $paras = $d.Paragraphs
    $count = $paras.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $paras.Item($i)
        if (($p.Range.Text -like "*$needle*") -and ($p.Range.Font.Size -eq $fontSize)) {
            return $p.Range
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) Extend the "From our top eight topics..." paragraph with a new run, and
#    insert the four new paragraphs that follow it (the Part 1 answer body).
# ---------------------------------------------------------------------------
$target1 = Find-ParagraphContaining("stands out to me the most")

$newBody1 = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">From our top eight topics we learned this class the one that stands out to me the most is ethical responsibility. </w:t></w:r><w:r><w:t xml:space="preserve"> I found ethical responsibility to be a prevalent topic during my time in this course. I had a very fuzzy understanding of this topic before this course. The activities and assignments greatly helped me better understand what ethical responsibility is.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:tab/><w:t xml:space="preserve">Ethical responsibility was always present throughout the term of this course. If it was not the main topic, it was always in the background or implied in other ways. The book covered this topic in detail dedicating a whole chapter to it. This reading was the main source of my understanding of ethical responsibility. Other than the book assignments related to the ACM code of ethics furthered my understanding. The </w:t></w:r><w:r><w:t>ACM code of ethics</w:t></w:r><w:r><w:t xml:space="preserve"> directly address the topic and provides guidelines for how to better yourself and others in the pursuit of professional ethics. I also found the discussion helpful general discussion covered this area, but I found the situational examples more helpful. </w:t></w:r><w:r><w:t>I enjoyed the small and large discussions and believe they helped me learn.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/></w:r></w:p>' + `
    '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t>Ethical responsibility is something I will undoubtedly encounter in the future, and this course has given me a greater understanding of what it means to be ethically responsibility.</w:t></w:r></w:p>'

Set-ParagraphXml "target1" $target1 $newBody1

# ---------------------------------------------------------------------------
# 2) The "Part 2" heading no longer starts a fresh rendered page (the page
#    break now falls earlier, inside the text we just added), so remove the
#    <w:lastRenderedPageBreak/> marker from that run.
# ---------------------------------------------------------------------------
$target2 = Find-ParagraphContainingWithFontSize("Essential question or workforce competencies", 20)

$newBody2 = '<w:p><w:pPr><w:keepNext/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>Part 2</w:t></w:r><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>: Essential question or workforce competencies</w:t></w:r></w:p>'

Set-ParagraphXml "target2" $target2 $newBody2

# ---------------------------------------------------------------------------
# 3) The paragraph beginning "This course " is now the one where the page
#    break actually renders, so it gains the <w:lastRenderedPageBreak/>.
# ---------------------------------------------------------------------------
$target3 = Find-ParagraphContaining("has offered many avenues to better understand")

$newBody3 = '<w:p><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="360"/><w:jc w:val="both"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">This course </w:t></w:r><w:r><w:t>has offered many avenues to better understand this essential question.</w:t></w:r><w:r><w:t xml:space="preserve"> For one the assigned readings seem to address this question to some extent in each of its sections. Discussion with the class also provided different viewpoints to this question that I had not considered myself. I personally found the scenario segment of the discussions to be the most useful. I enjoyed addressing specific case by case social and ethical issues, and I learned from my fellow students.</w:t></w:r></w:p>'

Set-ParagraphXml "target3" $target3 $newBody3

# ---------------------------------------------------------------------------
# 4) Likewise, the "Part 3" heading loses its now-stale page break marker.
# ---------------------------------------------------------------------------
$target4 = Find-ParagraphContainingWithFontSize("ACM Code of Ethics and Professional Conduct", 20)

$newBody4 = '<w:p><w:pPr><w:keepNext/><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve">Part </w:t></w:r><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve">3: </w:t></w:r><w:r><w:rPr><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>ACM Code of Ethics and Professional Conduct</w:t></w:r></w:p>'

Set-ParagraphXml "target4" $target4 $newBody4

Write-Host "Edit complete"
